# "Add files via upload" — populate the (previously empty) Sheet1 with a
# 9-row x 4-column block of the text "test" (A1:D9), matching the uploaded
# workbook, and restore the sheet's print/page setup (A4, portrait) that
# ships with that version of the file.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill A1:D9 with the literal string "test" (becomes a single shared string,
# referenced by all 36 cells).
$ws.Range("A1:D9").Value = "test"

# Page setup: A4 paper, portrait orientation.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Select the populated range, matching the saved selection in the workbook.
[void]$ws.Range("A1:D9").Select()
